$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = "Fer"
$ws.Range("B2").Value = "12.424.777-5"
$ws.Range("C2").Value = "CCTV/fotos/12424777-5/12.424.777-5_foto.jpg"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "123123"
# G2 keeps "Estudiante"

# Update row 3
$ws.Range("A3").Value = "Elon Musk"
$ws.Range("B3").Value = "12.123.123-1"
$ws.Range("C3").Value = "CCTV/fotos/12123123-1/12.123.123-1_foto.jpg"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "123123"
$ws.Range("E3").Font.Bold = $False
$ws.Range("F3").Font.Bold = $False
$ws.Range("G3").Value = "Docente"

# Remove old rows 4 and 5 entirely (no longer present in target)
$ws.Rows("4:5").Delete()
